$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before existing row 346 (shift the rest down, xlShiftDown = -4121)
$ws.Rows("346:348").Insert(-4121)

# Make sure the date column keeps the same date/time number format used elsewhere in column D
$ws.Range("D346:D348").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---- Row 346 : Hass / Especial, new week (2022-02-03), caja de 17 kilos ----
$ws.Range("A346").Value = 2
$ws.Range("B346").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C346").Value = "Coquimbo"
$ws.Range("D346").Value = 44595
$ws.Range("E346").Value = 4
$ws.Range("F346").Value = "Fruta"
$ws.Range("G346").Value = 100106
$ws.Range("H346").Value = "Oleaginosos"
$ws.Range("I346").Value = 100106002
$ws.Range("J346").Value = "Palta"
$ws.Range("K346").Value = "Hass"
$ws.Range("L346").Value = "Especial"
$ws.Range("M346").Value = 340
$ws.Range("N346").Value = 2500
$ws.Range("O346").Value = 2600
$ws.Range("P346").Value = 2550
$ws.Range("Q346").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R346").Value = "Provincia de Limarí"
$ws.Range("S346").Value = 2550
$ws.Range("T346").Value = 1

# ---- Row 347 : Hass / Primera, new week (2022-02-03), caja de 17 kilos ----
$ws.Range("A347").Value = 2
$ws.Range("B347").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C347").Value = "Coquimbo"
$ws.Range("D347").Value = 44595
$ws.Range("E347").Value = 4
$ws.Range("F347").Value = "Fruta"
$ws.Range("G347").Value = 100106
$ws.Range("H347").Value = "Oleaginosos"
$ws.Range("I347").Value = 100106002
$ws.Range("J347").Value = "Palta"
$ws.Range("K347").Value = "Hass"
$ws.Range("L347").Value = "Primera"
$ws.Range("M347").Value = 400
$ws.Range("N347").Value = 2200
$ws.Range("O347").Value = 2300
$ws.Range("P347").Value = 2250
$ws.Range("Q347").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R347").Value = "Provincia de Limarí"
$ws.Range("S347").Value = 2250
$ws.Range("T347").Value = 1

# ---- Row 348 : Hass / Segunda, new week (2022-02-03), caja de 17 kilos ----
$ws.Range("A348").Value = 2
$ws.Range("B348").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C348").Value = "Coquimbo"
$ws.Range("D348").Value = 44595
$ws.Range("E348").Value = 4
$ws.Range("F348").Value = "Fruta"
$ws.Range("G348").Value = 100106
$ws.Range("H348").Value = "Oleaginosos"
$ws.Range("I348").Value = 100106002
$ws.Range("J348").Value = "Palta"
$ws.Range("K348").Value = "Hass"
$ws.Range("L348").Value = "Segunda"
$ws.Range("M348").Value = 300
$ws.Range("N348").Value = 1800
$ws.Range("O348").Value = 1900
$ws.Range("P348").Value = 1850
$ws.Range("Q348").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R348").Value = "Provincia de Limarí"
$ws.Range("S348").Value = 1850
$ws.Range("T348").Value = 1
